$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 182 (pushes the existing rows 182:267 down
# to 183:268, and widens the sheet dimension from T267 to T268 automatically).
$ws.Rows("182:182").Insert()

# Populate the newly inserted row with the new weekly price-sheet entry.
$ws.Range("A182").Value = 4
$ws.Range("B182").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C182").Value = "Los Lagos"
$ws.Range("D182").Value = 44466
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = "Fruta"
$ws.Range("G182").Value = 100102
$ws.Range("H182").Value = "Cítricos"
$ws.Range("I182").Value = 100102003
$ws.Range("J182").Value = "Limón"
$ws.Range("K182").Value = "Sin especificar"
$ws.Range("L182").Value = "1a amarillo"
$ws.Range("M182").Value = 500
$ws.Range("N182").Value = 9000
$ws.Range("O182").Value = 9000
$ws.Range("P182").Value = 9000
$ws.Range("Q182").Value = '$/malla 18 kilos'
$ws.Range("R182").Value = "Provincia de Melipilla"
$ws.Range("S182").Value = 500
$ws.Range("T182").Value = 18
